# Update the si="28" shared string in cell B4 (Индикатор / "15.4.1 ...")
# to the corrected wording. The host engine automatically drops the
# now-unreferenced old shared string and appends the new text at the end
# of the shared-string table, renumbering every other <v> index that
# pointed past the removed slot — exactly matching the upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B4").Value = "15.4.1 Доля важных с точки зрения биологического разнообразия горных районов, находящихся под охраной"

# Move/record the active selection on the sheet to B4 (was B2).
$ws.Range("B4").Select()
